$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 170.33333
$ws.Range("I19").Value = 258.2
$ws.Range("J19").Value = 60.5
$ws.Range("K19").Value = 258.2
$ws.Range("L19").Value = 60.5
$ws.Range("M19").Value = -83.19999999999999
$ws.Range("N19").Value = -410.5
$ws.Range("H54").Value = 39999
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 39999
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 39999
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -40971
$ws.Range("H88").Value = 1848.9524
$ws.Range("I88").Value = 749
$ws.Range("J88").Value = 1964.7368
$ws.Range("K88").Value = 749
$ws.Range("L88").Value = 1964.7368
$ws.Range("M88").Value = -343
$ws.Range("N88").Value = -2776.7368
$ws.Range("H91").Value = 1848.9524
$ws.Range("I91").Value = 749
$ws.Range("J91").Value = 1964.7368
$ws.Range("K91").Value = 749
$ws.Range("L91").Value = 1964.7368
$ws.Range("M91").Value = 655
$ws.Range("N91").Value = -4772.7368
$ws.Range("H92").Value = 451
$ws.Range("I92").Value = 406.75
$ws.Range("K92").Value = 406.75
$ws.Range("M92").Value = 841.25
$ws.Range("H98").Value = 825.0769
$ws.Range("I98").Value = 798.2727
$ws.Range("J98").Value = 972.5
$ws.Range("K98").Value = 798.2727
$ws.Range("L98").Value = 972.5
$ws.Range("M98").Value = 699.7273
$ws.Range("N98").Value = -3968.5
$ws.Range("H122").Value = 825.0769
$ws.Range("I122").Value = 798.2727
$ws.Range("J122").Value = 972.5
$ws.Range("K122").Value = 2394.8181
$ws.Range("L122").Value = 2917.5
$ws.Range("M122").Value = 55.18190000000004
$ws.Range("N122").Value = -7817.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1013
$ws.Range("I2").Value = 370.2143
$ws.Range("K2").Value = 370.2143
$ws.Range("M2").Value = -257.2143
$ws.Range("H45").Value = 3499
$ws.Range("I45").Value = 1498.3334
$ws.Range("K45").Value = 1498.3334
$ws.Range("M45").Value = -1121.3334
$ws.Range("H88").Value = 5359.9
$ws.Range("I88").Value = 3126
$ws.Range("J88").Value = 6849.1665
$ws.Range("K88").Value = 3126
$ws.Range("L88").Value = 6849.1665
$ws.Range("M88").Value = -2720
$ws.Range("N88").Value = -7661.1665
$ws.Range("H91").Value = 5359.9
$ws.Range("I91").Value = 3126
$ws.Range("J91").Value = 6849.1665
$ws.Range("K91").Value = 3126
$ws.Range("L91").Value = 6849.1665
$ws.Range("M91").Value = -1722
$ws.Range("N91").Value = -9657.166499999999
$ws.Range("H97").Value = 10111.667
$ws.Range("I97").Value = 140
$ws.Range("K97").Value = 140
$ws.Range("M97").Value = 356
$ws.Range("H112").Value = 29789
$ws.Range("J112").Value = 29789
$ws.Range("L112").Value = 29789
$ws.Range("N112").Value = -32743
$ws.Range("H116").Value = 1013
$ws.Range("I116").Value = 370.2143
$ws.Range("K116").Value = 370.2143
$ws.Range("M116").Value = 1923.7857

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1013
$ws.Range("I3").Value = 370.2143
$ws.Range("K3").Value = 370.2143
$ws.Range("M3").Value = -256.2143
$ws.Range("H23").Value = 2972.3333
$ws.Range("J23").Value = 2972.3333
$ws.Range("L23").Value = 2972.3333
$ws.Range("N23").Value = -3538.3333
$ws.Range("H94").Value = 688.26666
$ws.Range("I94").Value = 722.4286
$ws.Range("K94").Value = 722.4286
$ws.Range("M94").Value = -271.4286
$ws.Range("H100").Value = 6412.5
$ws.Range("J100").Value = 6412.5
$ws.Range("L100").Value = 6412.5
$ws.Range("N100").Value = -8576.5
$ws.Range("H134").Value = 2357.7896
$ws.Range("I134").Value = 1674.875
$ws.Range("K134").Value = 5024.625
$ws.Range("M134").Value = -2489.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 566.3333
$ws.Range("I3").Value = 566.3333
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 566.3333
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -453.3333
$ws.Range("N3").ClearContents()
$ws.Range("H96").Value = 19833
$ws.Range("J96").Value = 19833
$ws.Range("L96").Value = 19833
$ws.Range("N96").Value = -25325

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 137.33333
$ws.Range("J63").Value = 150
$ws.Range("L63").Value = 450
$ws.Range("N63").Value = -1948
$ws.Range("H64").Value = 1575
$ws.Range("J64").Value = 1575
$ws.Range("L64").Value = 4725
$ws.Range("N64").Value = -5265
$ws.Range("H66").Value = 137.33333
$ws.Range("J66").Value = 150
$ws.Range("L66").Value = 1350
$ws.Range("N66").Value = -8838
$ws.Range("H67").Value = 1575
$ws.Range("J67").Value = 1575
$ws.Range("L67").Value = 4725
$ws.Range("N67").Value = -6597
$ws.Range("H107").Value = 1007.8889
$ws.Range("J107").Value = 1243.8462
$ws.Range("L107").Value = 3731.5386
$ws.Range("N107").Value = -7571.5386
$ws.Range("H113").Value = 255.66667
$ws.Range("J113").Value = 270
$ws.Range("L113").Value = 810
$ws.Range("N113").Value = -5150
$ws.Range("H131").Value = 1163.6428
$ws.Range("I131").Value = 1449
$ws.Range("J131").Value = 1116.0834
$ws.Range("K131").Value = 4347
$ws.Range("L131").Value = 3348.2502
$ws.Range("M131").Value = 693
$ws.Range("N131").Value = -13428.2502
$ws.Range("H137").Value = 4245
$ws.Range("I137").Value = 4245
$ws.Range("K137").Value = 12735
$ws.Range("M137").Value = -7635

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 231.5
$ws.Range("I13").Value = 100.5
$ws.Range("J13").Value = 362.5
$ws.Range("K13").Value = 100.5
$ws.Range("L13").Value = 362.5
$ws.Range("M13").Value = 38.5
$ws.Range("N13").Value = -640.5
$ws.Range("H22").Value = 2556.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1999
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H107").Value = 655.0769
$ws.Range("J107").Value = 1296.3334
$ws.Range("L107").Value = 1296.3334
$ws.Range("N107").Value = -5136.3334
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 2416.6875
$ws.Range("I122").Value = 2744.5386
$ws.Range("J122").Value = 996
$ws.Range("K122").Value = 8233.6158
$ws.Range("L122").Value = 2988
$ws.Range("M122").Value = -5783.6158
$ws.Range("N122").Value = -7888

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 825
$ws.Range("I22").Value = 825
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 825
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -530
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 825
$ws.Range("I27").Value = 825
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 825
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -718
$ws.Range("N27").ClearContents()
$ws.Range("H82").Value = 3881.8333
$ws.Range("I82").Value = 694
$ws.Range("J82").Value = 4519.4
$ws.Range("K82").Value = 694
$ws.Range("L82").Value = 4519.4
$ws.Range("M82").Value = -333
$ws.Range("N82").Value = -5241.4
$ws.Range("H85").Value = 3881.8333
$ws.Range("I85").Value = 694
$ws.Range("J85").Value = 4519.4
$ws.Range("K85").Value = 694
$ws.Range("L85").Value = 4519.4
$ws.Range("M85").Value = 554
$ws.Range("N85").Value = -7015.4
$ws.Range("H93").Value = 1379.0714
$ws.Range("I93").Value = 1301.125
$ws.Range("K93").Value = 1301.125
$ws.Range("M93").Value = -53.125
$ws.Range("H136").Value = 3313.625
$ws.Range("I136").Value = 2917.3333
$ws.Range("K136").Value = 8751.999899999999
$ws.Range("M136").Value = -6201.999899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 20000
$ws.Range("J4").Value = 20000
$ws.Range("L4").Value = 20000
$ws.Range("N4").Value = -20226
$ws.Range("H54").Value = 34999.062
$ws.Range("J54").Value = 34999.062
$ws.Range("L54").Value = 34999.062
$ws.Range("N54").Value = -36039.062
$ws.Range("H132").Value = 36976.605
$ws.Range("I132").Value = 42890.043
$ws.Range("J132").Value = 1496
$ws.Range("K132").Value = 128670.129
$ws.Range("L132").Value = 4488
$ws.Range("M132").Value = -126140.129
$ws.Range("N132").Value = -9548
$ws.Range("H136").Value = 2668.8
$ws.Range("I136").Value = 1386
$ws.Range("K136").Value = 4158
$ws.Range("M136").Value = -1608
